$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Progreso de tareas al 100% (antes 90%)
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 1

# La tarea de la fila 8 ya no esta "en proceso": ahora tiene 100 (completada)
$ws.Range("C8").Value = 100

# Nuevas tareas marcadas como "en proceso"
$ws.Range("C21").Value = "en proceso"
$ws.Range("C22").Value = "en proceso"

# Actualizar la celda seleccionada/activa de la hoja
$ws.Range("C9").Select()
